$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared by Overview!E2:F3 and the Status column on both language sheets)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(2,5).Value = $newStatus
$wsOverview.Cells.Item(2,6).Value = $newStatus
$wsOverview.Cells.Item(3,5).Value = $newStatus
$wsOverview.Cells.Item(3,6).Value = $newStatus
$wsZhCn.Cells.Item(2,3).Value = $newStatus
$wsZhCn.Cells.Item(3,3).Value = $newStatus
$wsDeDe.Cells.Item(2,3).Value = $newStatus
$wsDeDe.Cells.Item(3,3).Value = $newStatus

# ---------------------------------------------------------------------------
# 2) Column widths
#    Overview columns E/F and the "Status" column (C) on the language
#    sheets grow from ~17.2 to ~30; "Latest Target File"/"Latest Handback
#    File" columns (I/J) on the language sheets grow to 40.
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.15
    $ws.Columns.Item(9).ColumnWidth  = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
}

# ---------------------------------------------------------------------------
# 3) Fill in the "handback" columns (I = Latest Target File, J = Latest
#    Handback File, K = Latest Handback DateTime) for both language sheets,
#    and hyperlink the new Latest Target File entries.
# ---------------------------------------------------------------------------
$targetFileName = "e55b7646-a65e-4b63-b999-44d9629e963f.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fecafa7b044aa57c3563d19a74777e7efd8209c/e2e/e55b7646-a65e-4b63-b999-44d9629e963f.md"

# zh-cn
$wsZhCn.Cells.Item(2,9).Value = $targetFileName
$wsZhCn.Cells.Item(3,9).Value = $targetFileName
$wsZhCn.Cells.Item(2,10).Value = "e55b7646-a65e-4b63-b999-44d9629e963f.ddb342d63523bf2643e1ec407e955ae2e2a432df.zh-cn.xlf"
$wsZhCn.Cells.Item(3,10).Value = "e55b7646-a65e-4b63-b999-44d9629e963f.ddb342d63523bf2643e1ec407e955ae2e2a432df.zh-cn.xlf"
$wsZhCn.Cells.Item(2,11).Value = "2016-08-29 09:07:03"
$wsZhCn.Cells.Item(3,11).Value = "2016-08-29 09:07:03"

# de-de
$wsDeDe.Cells.Item(2,9).Value = $targetFileName
$wsDeDe.Cells.Item(3,9).Value = $targetFileName
$wsDeDe.Cells.Item(2,10).Value = "e55b7646-a65e-4b63-b999-44d9629e963f.ddb342d63523bf2643e1ec407e955ae2e2a432df.de-de.xlf"
$wsDeDe.Cells.Item(3,10).Value = "e55b7646-a65e-4b63-b999-44d9629e963f.ddb342d63523bf2643e1ec407e955ae2e2a432df.de-de.xlf"
$wsDeDe.Cells.Item(2,11).Value = "2016-08-29 09:07:15"
$wsDeDe.Cells.Item(3,11).Value = "2016-08-29 09:07:15"

# ---------------------------------------------------------------------------
# 4) Hyperlinks: rebuild each sheet's hyperlink list so it ends up in
#    A2, I2, A3, I3 order (matching relationship id order rId2..rId5),
#    adding new links on the "Latest Target File" cells.
# ---------------------------------------------------------------------------
foreach ($item in @(
        @{ ws = $wsZhCn; a2disp = "e55b7646-a65e-4b63-b999-44d9629e963f.md"; a2url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fecafa7b044aa57c3563d19a74777e7efd8209c/e2e/e55b7646-a65e-4b63-b999-44d9629e963f.md"; a3disp = "ffff5bbe20c6-f078-4066-aaea-8afb3fcc2bc7.md"; a3url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fecafa7b044aa57c3563d19a74777e7efd8209c/e2e/ffff5bbe20c6-f078-4066-aaea-8afb3fcc2bc7.md" },
        @{ ws = $wsDeDe; a2disp = "e55b7646-a65e-4b63-b999-44d9629e963f.md"; a2url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fecafa7b044aa57c3563d19a74777e7efd8209c/e2e/e55b7646-a65e-4b63-b999-44d9629e963f.md"; a3disp = "ffff5bbe20c6-f078-4066-aaea-8afb3fcc2bc7.md"; a3url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fecafa7b044aa57c3563d19a74777e7efd8209c/e2e/ffff5bbe20c6-f078-4066-aaea-8afb3fcc2bc7.md" }
    )) {
    $ws = $item.ws
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $item.a2url, "", "", $item.a2disp)
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetUrl, "", "", $targetFileName)
    $ws.Hyperlinks.Add($ws.Range("A3"), $item.a3url, "", "", $item.a3disp)
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetUrl, "", "", $targetFileName)
}
